$d = $word.ActiveDocument

# ---------------------------------------------------------------------------
# Step 0: remove the pre-existing "_GoBack" bookmark (it sits after "...as
# it's a different way of thinking." in the original document). It will be
# re-created later at its new location (inside the "immersive" run).
# ---------------------------------------------------------------------------
$existing = $d.Bookmarks("_GoBack")
$existing.Delete()

# ---------------------------------------------------------------------------
# Step 1: split "that the game they are making is immersive to the" into
# "...imme" | <bookmark _GoBack/> | "rsive to the"
# ---------------------------------------------------------------------------
$text = $d.Content.Text
$marker1 = "that the game they are making is imme"
$pos1 = $text.IndexOf($marker1) + $marker1.Length
$rng1 = $d.Range($pos1, $pos1)
$d.Bookmarks.Add("_GoBack", $rng1)

# ---------------------------------------------------------------------------
# Step 2: append a new sentence to the end of the "...large inheritance
# tress. " paragraph, as its own run.
# ---------------------------------------------------------------------------
$found2 = $d.Content.Find.Execute(
    "and large inheritance tress. ",
    $true, $false, $false, $false, $false, $true, 1, $false,
    "and large inheritance tress. In my personal opinion, ECS with FSM is not warranted, unless you have the capital to spend for 6 months or more development time to create an ECS paradigm with FSM coupled in and also, if you want an engine that has the capabilities of re-using different types of components together. ",
    2)

$text = $d.Content.Text
$marker2 = "and large inheritance tress. "
$pos2 = $text.IndexOf($marker2) + $marker2.Length
$rng2 = $d.Range($pos2, $pos2)
$d.Bookmarks.Add("ZZTempSplit", $rng2)
$d.Bookmarks("ZZTempSplit").Delete()

# ---------------------------------------------------------------------------
# Step 3: turn "...coupled together. The reason for this..." into
# "...coupled together" | " into their engine" | ". The reason for this..."
# ---------------------------------------------------------------------------
$found3 = $d.Content.Find.Execute(
    "with their own custom engine to integrate ECS with FSM coupled together. The reason",
    $true, $false, $false, $false, $false, $true, 1, $false,
    "with their own custom engine to integrate ECS with FSM coupled together into their engine. The reason",
    2)

$text = $d.Content.Text
$marker3a = "with their own custom engine to integrate ECS with FSM coupled together"
$pos3a = $text.IndexOf($marker3a) + $marker3a.Length
$rng3a = $d.Range($pos3a, $pos3a)
$d.Bookmarks.Add("ZZTempSplit", $rng3a)
$d.Bookmarks("ZZTempSplit").Delete()

$marker3b = " into their engine"
$pos3b = $pos3a + $marker3b.Length
$rng3b = $d.Range($pos3b, $pos3b)
$d.Bookmarks.Add("ZZTempSplit", $rng3b)
$d.Bookmarks("ZZTempSplit").Delete()

Write-Output "done"
